$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value2 = '30.417.61'
$ws.Cells.Item(2,5).Value2 = '  -1.85%  '
$ws.Cells.Item(3,4).Value2 = '1.910.51'
$ws.Cells.Item(3,5).Value2 = '  -2.42%  '
$ws.Cells.Item(4,5).Value2 = '  -0.05%  '
$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value2 = '238.79'
$c.Style = "Normal"
$ws.Cells.Item(5,5).Value2 = '  -2.39%  '
$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value2 = '1.001'
$c.Style = "Normal"
$ws.Cells.Item(6,5).Value2 = '  +0.03%  '
$c = $ws.Cells.Item(7,4)
$c.NumberFormat = "@"
$c.Value2 = '0.4761'
$c.Style = "Normal"
$ws.Cells.Item(7,5).Value2 = '  -2.27%  '
$c = $ws.Cells.Item(8,4)
$c.NumberFormat = "@"
$c.Value2 = '0.2853'
$c.Style = "Normal"
$ws.Cells.Item(8,5).Value2 = '  -3.42%  '
$c = $ws.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value2 = '0.06684'
$c.Style = "Normal"
$ws.Cells.Item(9,5).Value2 = '  -4.68%  '
$c = $ws.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value2 = '18.84'
$c.Style = "Normal"
$ws.Cells.Item(10,5).Value2 = '  -4.78%  '
$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value2 = '102.44'
$c.Style = "Normal"
$ws.Cells.Item(11,5).Value2 = '  -4.85%  '
$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value2 = '0.07729'
$c.Style = "Normal"
$ws.Cells.Item(12,5).Value2 = '  -1.03%  '
$ws.Cells.Item(13,4).Value2 = '1.915.23'
$ws.Cells.Item(13,5).Value2 = '  -2.24%  '
$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value2 = '5.206'
$c.Style = "Normal"
$ws.Cells.Item(14,5).Value2 = '  -5.11%  '
$c = $ws.Cells.Item(15,4)
$c.NumberFormat = "@"
$c.Value2 = '0.6737'
$c.Style = "Normal"
$ws.Cells.Item(15,5).Value2 = '  -3.98%  '
$ws.Cells.Item(16,4).Value2 = '30.419.54'
$ws.Cells.Item(16,5).Value2 = '  -1.88%  '
$c = $ws.Cells.Item(17,4)
$c.NumberFormat = "@"
$c.Value2 = '255.21'
$c.Style = "Normal"
$ws.Cells.Item(17,5).Value2 = '  -9.08%  '
$c = $ws.Cells.Item(18,4)
$c.NumberFormat = "@"
$c.Value2 = '1.000'
$c.Style = "Normal"
$ws.Cells.Item(18,5).Value2 = '  -0.04%  '
$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value2 = '0.000007496'
$c.Style = "Normal"
$ws.Cells.Item(19,5).Value2 = '  -3.95%  '
$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value2 = '12.72'
$c.Style = "Normal"
$ws.Cells.Item(20,5).Value2 = '  -4.48%  '
$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value2 = '5.395'
$c.Style = "Normal"
$ws.Cells.Item(21,5).Value2 = '  -2.96%  '
$ws.Cells.Item(22,5).Value2 = '  +0.03%  '
$ws.Cells.Item(23,2).Value2 = 'Chainlink'
$ws.Cells.Item(23,3).Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Cells.Item(23,4)
$c.NumberFormat = "@"
$c.Value2 = '6.302'
$c.Style = "Normal"
$ws.Cells.Item(23,5).Value2 = '  -3.37%  '
$ws.Cells.Item(24,2).Value2 = 'Cosmos'
$ws.Cells.Item(24,3).Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value2 = '9.449'
$c.Style = "Normal"
$ws.Cells.Item(24,5).Value2 = '  -4.01%  '
$ws.Cells.Item(25,2).Value2 = 'Monero'
$ws.Cells.Item(25,3).Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(25,4)
$c.NumberFormat = "@"
$c.Value2 = '164.96'
$c.Style = "Normal"
$ws.Cells.Item(25,5).Value2 = '  -2.02%  '
$ws.Cells.Item(26,2).Value2 = 'EthereumClassic'
$ws.Cells.Item(26,3).Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value2 = '18.93'
$c.Style = "Normal"
$ws.Cells.Item(26,5).Value2 = '  -5.10%  '
$ws.Cells.Item(27,2).Value2 = 'LidoDAOToken'
$ws.Cells.Item(27,3).Value2 = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value2 = '2.064'
$c.Style = "Normal"
$ws.Cells.Item(27,5).Value2 = '  -6.09%  '
$ws.Cells.Item(28,2).Value2 = 'Stellar'
$ws.Cells.Item(28,3).Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Cells.Item(28,4)
$c.NumberFormat = "@"
$c.Value2 = '0.1008'
$c.Style = "Normal"
$ws.Cells.Item(28,5).Value2 = '  -3.98%  '
$ws.Cells.Item(29,2).Value2 = 'Toncoin'
$ws.Cells.Item(29,3).Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Cells.Item(29,4)
$c.NumberFormat = "@"
$c.Value2 = '1.374'
$c.Style = "Normal"
$ws.Cells.Item(29,5).Value2 = '  -0.90%  '
$ws.Cells.Item(30,2).Value2 = 'Filecoin'
$ws.Cells.Item(30,3).Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Cells.Item(30,4)
$c.NumberFormat = "@"
$c.Value2 = '4.636'
$c.Style = "Normal"
$ws.Cells.Item(30,5).Value2 = '  -0.04%  '
$ws.Cells.Item(31,2).Value2 = 'PancakeSwap'
$ws.Cells.Item(31,3).Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Cells.Item(31,4)
$c.NumberFormat = "@"
$c.Value2 = '1.511'
$c.Style = "Normal"
$ws.Cells.Item(31,5).Value2 = '  -3.70%  '
$ws.Cells.Item(32,2).Value2 = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32,3).Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(32,4)
$c.NumberFormat = "@"
$c.Value2 = '4.265'
$c.Style = "Normal"
$ws.Cells.Item(32,5).Value2 = '  -4.64%  '
$ws.Cells.Item(33,2).Value2 = 'Hedera'
$ws.Cells.Item(33,3).Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(33,4)
$c.NumberFormat = "@"
$c.Value2 = '0.04783'
$c.Style = "Normal"
$ws.Cells.Item(33,5).Value2 = '  -2.93%  '
$ws.Cells.Item(34,2).Value2 = 'ImmutableX'
$ws.Cells.Item(34,3).Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Cells.Item(34,4)
$c.NumberFormat = "@"
$c.Value2 = '0.7316'
$c.Style = "Normal"
$ws.Cells.Item(34,5).Value2 = '  -3.25%  '
$ws.Cells.Item(35,2).Value2 = 'ARBITRUM'
$ws.Cells.Item(35,3).Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value2 = '1.113'
$c.Style = "Normal"
$ws.Cells.Item(35,5).Value2 = '  -5.09%  '
$ws.Cells.Item(36,2).Value2 = 'Frax'
$ws.Cells.Item(36,3).Value2 = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Cells.Item(36,4)
$c.NumberFormat = "@"
$c.Value2 = '0.9995'
$c.Style = "Normal"
$ws.Cells.Item(36,5).Value2 = '  -0.08%  '
$ws.Cells.Item(37,2).Value2 = 'HuobiToken'
$ws.Cells.Item(37,3).Value2 = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Cells.Item(37,4)
$c.NumberFormat = "@"
$c.Value2 = '2.712'
$c.Style = "Normal"
$ws.Cells.Item(37,5).Value2 = '  -0.80%  '
$ws.Cells.Item(38,2).Value2 = 'VeChain'
$ws.Cells.Item(38,3).Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(38,4)
$c.NumberFormat = "@"
$c.Value2 = '0.01928'
$c.Style = "Normal"
$ws.Cells.Item(38,5).Value2 = '  -4.18%  '
$ws.Cells.Item(39,2).Value2 = 'MXToken'
$ws.Cells.Item(39,3).Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Cells.Item(39,4)
$c.NumberFormat = "@"
$c.Value2 = '2.598'
$c.Style = "Normal"
$ws.Cells.Item(39,5).Value2 = '  -3.48%  '
$ws.Cells.Item(40,2).Value2 = 'FraxShare'
$ws.Cells.Item(40,3).Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value2 = '6.250'
$c.Style = "Normal"
$ws.Cells.Item(40,5).Value2 = '  -4.62%  '
$ws.Cells.Item(41,2).Value2 = 'Aave'
$ws.Cells.Item(41,3).Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Cells.Item(41,4)
$c.NumberFormat = "@"
$c.Value2 = '74.38'
$c.Style = "Normal"
$ws.Cells.Item(41,5).Value2 = '  -4.89%  '
$ws.Cells.Item(42,2).Value2 = 'RenderToken'
$ws.Cells.Item(42,3).Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(42,4)
$c.NumberFormat = "@"
$c.Value2 = '1.997'
$c.Style = "Normal"
$ws.Cells.Item(42,5).Value2 = '  -7.09%  '
$ws.Cells.Item(43,2).Value2 = 'TrustWalletToken'
$ws.Cells.Item(43,3).Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value2 = '0.8649'
$c.Style = "Normal"
$ws.Cells.Item(43,5).Value2 = '  -4.20%  '
$ws.Cells.Item(44,2).Value2 = 'Quant'
$ws.Cells.Item(44,3).Value2 = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value2 = '106.68'
$c.Style = "Normal"
$ws.Cells.Item(44,5).Value2 = '  -2.45%  '
$ws.Cells.Item(45,2).Value2 = 'Maker'
$ws.Cells.Item(45,3).Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(45,4).Value2 = '1.061.68'
$ws.Cells.Item(45,5).Value2 = '  +5.18%  '
$ws.Cells.Item(46,2).Value2 = 'TheSandbox'
$ws.Cells.Item(46,3).Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Cells.Item(46,4)
$c.NumberFormat = "@"
$c.Value2 = '0.4250'
$c.Style = "Normal"
$ws.Cells.Item(46,5).Value2 = '  -4.86%  '
$ws.Cells.Item(47,2).Value2 = 'PaxDollar'
$ws.Cells.Item(47,3).Value2 = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Cells.Item(47,4)
$c.NumberFormat = "@"
$c.Value2 = '1.000'
$c.Style = "Normal"
$ws.Cells.Item(47,5).Value2 = '  -0.04%  '
$ws.Cells.Item(48,2).Value2 = 'Aptos'
$ws.Cells.Item(48,3).Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value2 = '7.488'
$c.Style = "Normal"
$ws.Cells.Item(48,5).Value2 = '  -8.09%  '
$ws.Cells.Item(49,2).Value2 = 'Algorand'
$ws.Cells.Item(49,3).Value2 = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Cells.Item(49,4)
$c.NumberFormat = "@"
$c.Value2 = '0.1198'
$c.Style = "Normal"
$ws.Cells.Item(49,5).Value2 = '  -4.51%  '
$ws.Cells.Item(50,2).Value2 = 'Elrond'
$ws.Cells.Item(50,3).Value2 = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Cells.Item(50,4)
$c.NumberFormat = "@"
$c.Value2 = '34.92'
$c.Style = "Normal"
$ws.Cells.Item(50,5).Value2 = '  -2.85%  '
$ws.Cells.Item(51,2).Value2 = 'EnergySwap'
$ws.Cells.Item(51,3).Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value2 = '8.855'
$c.Style = "Normal"
$ws.Cells.Item(51,5).Value2 = '  -5.02%  '
